$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "虚无" (Nothingness) buff row (ID 0) ---
$ws.Rows(2).Delete()

# --- Remove the "减伤" (damage-reduction ChangeProperty) row and the
#     "减防" (defense-reduction ChangeProperty) row. After the first
#     deletion they sit at rows 10 and 8 respectively; delete the lower
#     one first so the second delete's row index stays correct. ---
$ws.Rows(10).Delete()
$ws.Rows(8).Delete()

# --- Renumber the ID column (A) for every remaining data row (2-12) so
#     the sequence is contiguous again (0..10) after the three deletes. ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10

# --- Append the five new buffs as new rows 13-17 ---
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "力竭"
$ws.Range("C13").Value = "Exhaustion"
$ws.Range("D13").Value = -1
$ws.Range("E13").Value = -1
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "buff_confine_icon"
$ws.Range("K13").Value = "力竭"

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "虚弱"
$ws.Range("C14").Value = "Weakness"
$ws.Range("D14").Value = -1
$ws.Range("E14").Value = -1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "buff_confine_icon"
$ws.Range("K14").Value = "虚弱"

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "破绽"
$ws.Range("C15").Value = "Flaw"
$ws.Range("D15").Value = -1
$ws.Range("E15").Value = -1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "buff_confine_icon"
$ws.Range("K15").Value = "破绽"

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "连击"
$ws.Range("C16").Value = "Batter"
$ws.Range("D16").Value = -1
$ws.Range("E16").Value = -1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "buff_confine_icon"
$ws.Range("K16").Value = "反击"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "锋锐"
$ws.Range("C17").Value = "Sharp"
$ws.Range("D17").Value = -1
$ws.Range("E17").Value = -1
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "buff_confine_icon"
$ws.Range("K17").Value = "锋锐"

# --- Update the active selection to match the final edit location ---
$ws.Range("G17").Select() | Out-Null
